$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "edit-fields": insert 4 new rows before row 24
#   24 alert    / node_id  / preset / 1
#   25 barcode  / node_id  / type   / string
#   26 barcode  / model_id / type   / string
#   27 barcode  / barcode  / type   / string
# (existing rows 24-39 shift down to 28-43)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("edit-fields")
$ws2.Range("A24:D27").EntireRow.Insert() | Out-Null

$ws2.Range("A24").Value = "alert"
$ws2.Range("B24").Value = "node_id"
$ws2.Range("C24").Value = "preset"
$ws2.Range("D24").Value = 1

# ---------------------------------------------------------------
# Sheet "extras": insert 1 new row after row 1 (new row 2)
#   2 alert-action / content / class / textarea
# (existing rows 2-38 shift down to 3-39)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("extras")
$ws3.Range("A2:D2").EntireRow.Insert() | Out-Null

# Re-apply the formatting the inserted row should have (matching the row
# that used to be row 2 and is now row 4), since Excel's auto-insert
# picks up formatting from the row above (the header) by default.
$ws3.Range("A4:D4").Copy() | Out-Null
$ws3.Range("A2:D2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws3.Range("A2").Value = "alert-action"
$ws3.Range("B2").Value = "content"
$ws3.Range("C2").Value = "class"
$ws3.Range("D2").Value = "textarea"

# back to edit-fields to finish the barcode rows
$ws2.Range("B26").Value = "model_id"

$ws2.Range("A25").Value = "barcode"
$ws2.Range("B25").Value = "node_id"
$ws2.Range("C25").Value = "type"
$ws2.Range("D25").Value = "string"

$ws2.Range("A26").Value = "barcode"
$ws2.Range("C26").Value = "type"
$ws2.Range("D26").Value = "string"

$ws2.Range("A27").Value = "barcode"
$ws2.Range("B27").Value = "barcode"
$ws2.Range("C27").Value = "type"
$ws2.Range("D27").Value = "string"

$ws2.Range("D28").Select() | Out-Null
$ws3.Range("B2").Select() | Out-Null
